# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Updates the "Estado de Cuenta" worker/period table (rows 16-22, columns
# B:G) on Hoja1: the existing rows for EDUARDO LUIS RIOS ARRIETA / ROSNI DEL
# CARMEN PINTO SARMIENTO are re-sequenced by period (1801..1805) and a new
# row is inserted for ROSNI's period 1801, with her "Valor Mora" updated
# from 750000 to 781242 across all of her periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: CC 73116535 - EDUARDO LUIS RIOS ARRIETA - periodo 1801
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73116535"
$ws.Range("D16").Value = "EDUARDO LUIS RIOS ARRIETA"
$ws.Range("E16").Value = "1801"
$ws.Range("F16").Value = 29509
$ws.Range("G16").Value = 737717

# Row 17: CC 45750363 - ROSNI DEL CARMEN PINTO SARMIENTO - periodo 1801 (new)
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45750363"
$ws.Range("D17").Value = "ROSNI DEL CARMEN PINTO SARMIENTO"
$ws.Range("E17").Value = "1801"
$ws.Range("F17").Value = 30000
$ws.Range("G17").Value = 781242

# Row 18: CC 73116535 - EDUARDO LUIS RIOS ARRIETA - periodo 1802
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73116535"
$ws.Range("D18").Value = "EDUARDO LUIS RIOS ARRIETA"
$ws.Range("E18").Value = "1802"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 737717

# Row 19: CC 45750363 - ROSNI DEL CARMEN PINTO SARMIENTO - periodo 1802
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45750363"
$ws.Range("D19").Value = "ROSNI DEL CARMEN PINTO SARMIENTO"
$ws.Range("E19").Value = "1802"
$ws.Range("F19").Value = 30000
$ws.Range("G19").Value = 781242

# Row 20: CC 45750363 - ROSNI DEL CARMEN PINTO SARMIENTO - periodo 1803
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "45750363"
$ws.Range("D20").Value = "ROSNI DEL CARMEN PINTO SARMIENTO"
$ws.Range("E20").Value = "1803"
$ws.Range("F20").Value = 30000
$ws.Range("G20").Value = 781242

# Row 21: CC 45750363 - ROSNI DEL CARMEN PINTO SARMIENTO - periodo 1804
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "45750363"
$ws.Range("D21").Value = "ROSNI DEL CARMEN PINTO SARMIENTO"
$ws.Range("E21").Value = "1804"
$ws.Range("F21").Value = 30000
$ws.Range("G21").Value = 781242

# Row 22: CC 45750363 - ROSNI DEL CARMEN PINTO SARMIENTO - periodo 1805
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "45750363"
$ws.Range("D22").Value = "ROSNI DEL CARMEN PINTO SARMIENTO"
$ws.Range("E22").Value = "1805"
$ws.Range("F22").Value = 30000
$ws.Range("G22").Value = 781242
